# Adapt column header formatting to respective input file names.
# - Rename "..._old" headers (cols A-J) to "..._FV2404"
# - Rename "..._new" headers (cols L-U) to "..._FV2410"
# - Wrap the used range in an Excel Table (Table1)
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
    "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
    "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the data range into an Excel Table (ListObject)
$dataRange = $ws.Range("A1:U68")
$tbl = $ws.ListObjects.Add(1, $dataRange, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
$tbl.ShowTableStyleRowStripes = $true
$tbl.ShowTableStyleColumnStripes = $false
$tbl.ShowTableStyleFirstColumn = $false
$tbl.ShowTableStyleLastColumn = $false

# Freeze the top (header) row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
